# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.745.02"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "3.798.44"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'704.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "'169.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("D7").Value = "3.795.78"
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").Value = "'7.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").Value = "'36.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "4.438.84"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "3.810.94"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "70.718.24"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").Value = "'17.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").Value = "'497.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("E22").Value = "  -4.71%  "
$ws.Range("D23").Value = "'0.727"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'84.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").Value = "'10.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").Value = "3.949.10"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'2.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("E32").Value = "  -4.24%  "
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").Value = "'0.172"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.70%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "3.768.58"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "'9.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("D39").Value = "'0.102"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "'1.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").Value = "'3.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("B44").Value = "FLOKI"
$ws.Range("C44").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D44").Value = "'0.000325"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'164.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "'425.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").Value = "'48.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'8.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("E51").Value = "  -1.67%  "

Write-Host "Applied crypto list updates"